# Apply updated cryptocurrency price/volume figures (cryptos.xlsx refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a cell as literal text. A leading apostrophe forces Excel to
# keep values that look numeric (e.g. "0.4586") stored as text, matching the
# original inline-string cells instead of being auto-converted to numbers.
function Set-TextCell([string]$addr, [string]$text) {
    $ws.Range($addr).Value = "'" + $text
}

$ws.Range('D2').Value = '28.887.06'
$ws.Range('E2').Value = '  -1.31%  '
$ws.Range('D3').Value = '1.907.34'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  +0.15%  '
Set-TextCell 'D5' '324.71'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('E6').Value = '  +0.17%  '
Set-TextCell 'D7' '0.4586'
$ws.Range('E7').Value = '  -0.69%  '
$ws.Range('E8').Value = '  -1.41%  '
Set-TextCell 'D9' '0.07722'
$ws.Range('E9').Value = '  -1.04%  '
Set-TextCell 'D10' '0.9801'
$ws.Range('E10').Value = '  +0.67%  '
Set-TextCell 'D11' '22.17'
$ws.Range('E11').Value = '  -1.86%  '
$ws.Range('D12').Value = '1.873.95'
$ws.Range('E12').Value = '  -3.23%  '
Set-TextCell 'D13' '5.678'
$ws.Range('E13').Value = '  -1.76%  '
Set-TextCell 'D14' '6.952'
$ws.Range('E14').Value = '  -1.63%  '
Set-TextCell 'D15' '0.07059'
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('E16').Value = '  +0.14%  '
Set-TextCell 'D17' '83.75'
Set-TextCell 'D18' '0.000009469'
$ws.Range('E18').Value = '  -2.39%  '
Set-TextCell 'D19' '16.63'
$ws.Range('E19').Value = '  -2.27%  '
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = '28.887.48'
$ws.Range('E21').Value = '  -1.39%  '
Set-TextCell 'D22' '5.321'
$ws.Range('E22').Value = '  -2.87%  '
Set-TextCell 'D23' '10.92'
$ws.Range('E23').Value = '  -1.14%  '
Set-TextCell 'D24' '2.097'
$ws.Range('E24').Value = '  +0.31%  '
Set-TextCell 'D25' '158.78'
$ws.Range('E25').Value = '  +1.04%  '
Set-TextCell 'D26' '19.01'
$ws.Range('E26').Value = '  -1.73%  '
Set-TextCell 'D27' '5.669'
$ws.Range('E27').Value = '  -1.67%  '
Set-TextCell 'D28' '117.77'
$ws.Range('E28').Value = '  -0.91%  '
Set-TextCell 'D29' '1.869'
$ws.Range('E29').Value = '  +1.30%  '
Set-TextCell 'D30' '0.09301'
$ws.Range('E30').Value = '  -0.43%  '
Set-TextCell 'D31' '0.8642'
$ws.Range('E31').Value = '  +0.25%  '
Set-TextCell 'D32' '5.088'
$ws.Range('E32').Value = '  -1.42%  '
Set-TextCell 'D33' '1.246'
$ws.Range('E33').Value = '  -4.18%  '
Set-TextCell 'D34' '3.015'
$ws.Range('E34').Value = '  -1.90%  '
Set-TextCell 'D35' '0.05712'
$ws.Range('E35').Value = '  -0.98%  '
Set-TextCell 'D36' '1.156'
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('E37').Value = '  +0.16%  '
Set-TextCell 'D38' '0.02043'
$ws.Range('E38').Value = '  -1.66%  '
Set-TextCell 'D39' '7.452'
$ws.Range('E39').Value = '  -2.39%  '
Set-TextCell 'D40' '0.5491'
$ws.Range('E40').Value = '  -2.70%  '
Set-TextCell 'D41' '0.1753'
$ws.Range('E41').Value = '  -1.33%  '
Set-TextCell 'D42' '2.865'
$ws.Range('E42').Value = '  +5.83%  '
Set-TextCell 'D43' '9.341'
$ws.Range('E43').Value = '  -0.50%  '
Set-TextCell 'D44' '0.000002795'
$ws.Range('E44').Value = '  -10.00%  '
Set-TextCell 'D45' '2.165'
$ws.Range('E45').Value = '  +4.40%  '
Set-TextCell 'D46' '0.5163'
$ws.Range('E46').Value = '  -1.88%  '
Set-TextCell 'D47' '11.25'
$ws.Range('E47').Value = '  -1.60%  '
Set-TextCell 'D48' '0.06890'
$ws.Range('E48').Value = '  +0.44%  '
Set-TextCell 'D49' '110.47'
$ws.Range('E49').Value = '  -0.63%  '
Set-TextCell 'D50' '1.777'
$ws.Range('E50').Value = '  -1.90%  '
Set-TextCell 'D51' '0.2865'
$ws.Range('E51').Value = '  -4.35%  '
